$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "type" column (C) values: varChar -> nvarChar, plus a couple of
# null/width changes, per the authored edit.

# Users table
$ws.Range("C5").Value = "nvarChar(50), non-null"
$ws.Range("C6").Value = "nvarChar(50), non-null"
$ws.Range("C7").Value = "nvarChar(50), non-null"
$ws.Range("C8").Value = "nvarChar(50), non-null"

# Contacts table
$ws.Range("C18").Value = "nvarChar(50), non-null"
$ws.Range("C19").Value = "nvarChar(50), non-null"
$ws.Range("C20").Value = "nvarChar(10), non-null"
$ws.Range("C21").Value = "nvarChar(10), null"
$ws.Range("C22").Value = "nvarChar(10), null"
$ws.Range("C23").Value = "nvarChar(100), non-null"
$ws.Range("C24").Value = "nvarChar(100), non-null"
$ws.Range("C25").Value = "nvarChar(5), non-null"
$ws.Range("C26").Value = "nvarChar(100), non-null"

# Update selected cell to match author's saved view state
$ws.Range("E11").Select()
